$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 14).Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 7338.8
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 13).Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2648
$ws.Cells.Item(70, 9).Value = 2135.3
$ws.Cells.Item(70, 11).Value = 6405.900000000001
$ws.Cells.Item(70, 13).Value = -6135.900000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 2648
$ws.Cells.Item(73, 9).Value = 2135.3
$ws.Cells.Item(73, 11).Value = 6405.900000000001
$ws.Cells.Item(73, 13).Value = -5469.900000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 8999.6
$ws.Cells.Item(86, 9).Value = 8333
$ws.Cells.Item(86, 10).Value = 9999.5
$ws.Cells.Item(86, 11).Value = 8333
$ws.Cells.Item(86, 12).Value = 9999.5
$ws.Cells.Item(86, 13).Value = -7210
$ws.Cells.Item(86, 14).Value = -12245.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 8999.6
$ws.Cells.Item(89, 9).Value = 8333
$ws.Cells.Item(89, 10).Value = 9999.5
$ws.Cells.Item(89, 11).Value = 41665
$ws.Cells.Item(89, 12).Value = 49997.5
$ws.Cells.Item(89, 13).Value = -36049
$ws.Cells.Item(89, 14).Value = -61229.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 3710.318
$ws.Cells.Item(112, 10).Value = 3696.524
$ws.Cells.Item(112, 12).Value = 11089.572
$ws.Cells.Item(112, 14).Value = -13305.572

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 4807.5
$ws.Cells.Item(125, 9).Value = 3880.9092
$ws.Cells.Item(125, 11).Value = 34928.1828
$ws.Cells.Item(125, 13).Value = -32468.1828

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2714.3262
$ws.Cells.Item(138, 9).Value = 1505.375
$ws.Cells.Item(138, 10).Value = 4033.182
$ws.Cells.Item(138, 11).Value = 4516.125
$ws.Cells.Item(138, 12).Value = 12099.546
$ws.Cells.Item(138, 13).Value = 623.875
$ws.Cells.Item(138, 14).Value = -22379.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 2858.1428
$ws.Cells.Item(13, 9).Value = 1168
$ws.Cells.Item(13, 10).Value = 12999
$ws.Cells.Item(13, 11).Value = 1168
$ws.Cells.Item(13, 12).Value = 12999
$ws.Cells.Item(13, 13).Value = -1024
$ws.Cells.Item(13, 14).Value = -13287

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 14).Value = ""
$ws.Cells.Item(25, 12).Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 3151.3027
$ws.Cells.Item(74, 9).Value = 2482.1833
$ws.Cells.Item(74, 11).Value = 2482.1833
$ws.Cells.Item(74, 13).Value = -1608.1833

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 54000
$ws.Cells.Item(76, 10).Value = 54000
$ws.Cells.Item(76, 12).Value = 54000
$ws.Cells.Item(76, 14).Value = -54676

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 3151.3027
$ws.Cells.Item(77, 9).Value = 2482.1833
$ws.Cells.Item(77, 11).Value = 12410.9165
$ws.Cells.Item(77, 13).Value = -8042.916500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(79, 8).Value = 54000
$ws.Cells.Item(79, 10).Value = 54000
$ws.Cells.Item(79, 12).Value = 54000
$ws.Cells.Item(79, 14).Value = -56340

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 5660.1724
$ws.Cells.Item(110, 9).Value = 4903.7
$ws.Cells.Item(110, 10).Value = 7341.222
$ws.Cells.Item(110, 11).Value = 4903.7
$ws.Cells.Item(110, 12).Value = 7341.222
$ws.Cells.Item(110, 13).Value = -2858.7
$ws.Cells.Item(110, 14).Value = -11431.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3065.7937
$ws.Cells.Item(132, 9).Value = 2994.94
$ws.Cells.Item(132, 11).Value = 8984.82
$ws.Cells.Item(132, 13).Value = -6454.82

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 10271.286
$ws.Cells.Item(37, 9).Value = 2974
$ws.Cells.Item(37, 11).Value = 2974
$ws.Cells.Item(37, 13).Value = -2837

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1357.3448
$ws.Cells.Item(94, 9).Value = 1126.1569
$ws.Cells.Item(94, 10).Value = 3041.7144
$ws.Cells.Item(94, 11).Value = 1126.1569
$ws.Cells.Item(94, 12).Value = 3041.7144
$ws.Cells.Item(94, 13).Value = -675.1569
$ws.Cells.Item(94, 14).Value = -3943.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 1875.25
$ws.Cells.Item(4, 9).Value = 1667
$ws.Cells.Item(4, 11).Value = 1667
$ws.Cells.Item(4, 13).Value = -1555

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 64962.5
$ws.Cells.Item(23, 10).Value = 64954.332
$ws.Cells.Item(23, 12).Value = 64954.332
$ws.Cells.Item(23, 14).Value = -65434.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(27, 8).Value = 64962.5
$ws.Cells.Item(27, 10).Value = 64954.332
$ws.Cells.Item(27, 12).Value = 64954.332
$ws.Cells.Item(27, 14).Value = -65338.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 5837.926
$ws.Cells.Item(58, 9).Value = 4244.5454
$ws.Cells.Item(58, 10).Value = 6933.375
$ws.Cells.Item(58, 11).Value = 4244.5454
$ws.Cells.Item(58, 12).Value = 6933.375
$ws.Cells.Item(58, 13).Value = -4041.5454
$ws.Cells.Item(58, 14).Value = -7339.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 2480.818
$ws.Cells.Item(105, 9).Value = 2228.9
$ws.Cells.Item(105, 11).Value = 2228.9
$ws.Cells.Item(105, 13).Value = -481.9000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 4018
$ws.Cells.Item(107, 9).Value = 5277
$ws.Cells.Item(107, 10).Value = 1500
$ws.Cells.Item(107, 11).Value = 5277
$ws.Cells.Item(107, 12).Value = 1500
$ws.Cells.Item(107, 13).Value = -3357
$ws.Cells.Item(107, 14).Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 5837.926
$ws.Cells.Item(136, 9).Value = 4244.5454
$ws.Cells.Item(136, 10).Value = 6933.375
$ws.Cells.Item(136, 11).Value = 12733.6362
$ws.Cells.Item(136, 12).Value = 20800.125
$ws.Cells.Item(136, 13).Value = -10183.6362
$ws.Cells.Item(136, 14).Value = -25900.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 1999
$ws.Cells.Item(20, 9).Value = 1998
$ws.Cells.Item(20, 11).Value = 5994
$ws.Cells.Item(20, 13).Value = -5767

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(21, 8).Value = 449
$ws.Cells.Item(21, 9).Value = 449
$ws.Cells.Item(21, 11).Value = 1347
$ws.Cells.Item(21, 13).Value = -1174

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 286100.16
$ws.Cells.Item(26, 10).Value = 573.75
$ws.Cells.Item(26, 12).Value = 1721.25
$ws.Cells.Item(26, 14).Value = -2297.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 2198
$ws.Cells.Item(113, 10).Value = 3369.3333
$ws.Cells.Item(113, 12).Value = 10107.9999
$ws.Cells.Item(113, 14).Value = -14447.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 176920980
$ws.Cells.Item(11, 9).Value = 212105180
$ws.Cells.Item(11, 10).Value = 1000000
$ws.Cells.Item(11, 11).Value = 212105180
$ws.Cells.Item(11, 12).Value = 1000000
$ws.Cells.Item(11, 13).Value = -212105041
$ws.Cells.Item(11, 14).Value = -1000278

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 14).Value = ""
$ws.Cells.Item(18, 12).Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3305.0625
$ws.Cells.Item(113, 9).Value = 2715.9
$ws.Cells.Item(113, 11).Value = 2715.9
$ws.Cells.Item(113, 13).Value = -545.9000000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1991.4783
$ws.Cells.Item(132, 9).Value = 1886.9722
$ws.Cells.Item(132, 11).Value = 5660.9166
$ws.Cells.Item(132, 13).Value = -3130.9166

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 253750.25
$ws.Cells.Item(61, 9).Value = 336667.66
$ws.Cells.Item(61, 10).Value = 4998
$ws.Cells.Item(61, 11).Value = 336667.66
$ws.Cells.Item(61, 12).Value = 4998
$ws.Cells.Item(61, 13).Value = -336465.66
$ws.Cells.Item(61, 14).Value = -5402

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1937.8334
$ws.Cells.Item(68, 9).Value = 1902.125
$ws.Cells.Item(68, 11).Value = 1902.125
$ws.Cells.Item(68, 13).Value = -1153.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 1937.8334
$ws.Cells.Item(71, 9).Value = 1902.125
$ws.Cells.Item(71, 11).Value = 9510.625
$ws.Cells.Item(71, 13).Value = -5766.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 253750.25
$ws.Cells.Item(113, 9).Value = 336667.66
$ws.Cells.Item(113, 10).Value = 4998
$ws.Cells.Item(113, 11).Value = 336667.66
$ws.Cells.Item(113, 12).Value = 4998
$ws.Cells.Item(113, 13).Value = -334497.66
$ws.Cells.Item(113, 14).Value = -9338

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(119, 8).Value = 54200
$ws.Cells.Item(119, 10).Value = 54200
$ws.Cells.Item(119, 12).Value = 54200
$ws.Cells.Item(119, 14).Value = -63876

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127, 8).Value = 60000
$ws.Cells.Item(127, 10).Value = 60000
$ws.Cells.Item(127, 12).Value = 60000
$ws.Cells.Item(127, 14).Value = -69920

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 9249.960999999999
$ws.Cells.Item(132, 9).Value = 8640.643
$ws.Cells.Item(132, 11).Value = 25921.929
$ws.Cells.Item(132, 13).Value = -23391.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 8859.583000000001
$ws.Cells.Item(136, 9).Value = 8859.583000000001
$ws.Cells.Item(136, 11).Value = 26578.749
$ws.Cells.Item(136, 13).Value = -24028.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 29999
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 1000
$ws.Cells.Item(18, 9).Value = 1000
$ws.Cells.Item(18, 11).Value = 1000
$ws.Cells.Item(18, 13).Value = -827

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(34, 8).Value = 40000
$ws.Cells.Item(34, 9).Value = 40000
$ws.Cells.Item(34, 11).Value = 40000
$ws.Cells.Item(34, 13).Value = -39797

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4623.0586
$ws.Cells.Item(62, 10).Value = 4397.222
$ws.Cells.Item(62, 12).Value = 4397.222
$ws.Cells.Item(62, 14).Value = -5645.222

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 4623.0586
$ws.Cells.Item(65, 10).Value = 4397.222
$ws.Cells.Item(65, 12).Value = 21986.11
$ws.Cells.Item(65, 14).Value = -28226.11

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 875.44446
$ws.Cells.Item(113, 9).Value = 973
$ws.Cells.Item(113, 10).Value = 722.1429000000001
$ws.Cells.Item(113, 11).Value = 2919
$ws.Cells.Item(113, 12).Value = 2166.4287
$ws.Cells.Item(113, 13).Value = -749
$ws.Cells.Item(113, 14).Value = -6506.4287

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3016.2
$ws.Cells.Item(132, 9).Value = 3428.5151
$ws.Cells.Item(132, 10).Value = 1882.3334
$ws.Cells.Item(132, 11).Value = 10285.5453
$ws.Cells.Item(132, 12).Value = 5647.0002
$ws.Cells.Item(132, 13).Value = -7755.5453
$ws.Cells.Item(132, 14).Value = -10707.0002

